# Removing less than USD 5 price from extrapolation calibration because it is just a noise
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9
$ws.Range("D9").Value = 119331.1888996099
$ws.Range("E9").Value = -0.08567331759433981
$ws.Range("F9").Value = 0.3360058191039166
$ws.Range("G9").Value = -1.557671579805229
$ws.Range("H9").Value = 10.04200035094302

# Row 10
$ws.Range("D10").Value = 120823.8872333461
$ws.Range("E10").Value = -0.1221290607987758
$ws.Range("F10").Value = 0.4423840516989567
$ws.Range("G10").Value = -1.883259025656171
$ws.Range("H10").Value = 9.447172317139584

# Row 11
$ws.Range("D11").Value = 122766.0462835769
$ws.Range("E11").Value = -0.1971808894111483
$ws.Range("F11").Value = 0.765067017845186
$ws.Range("G11").Value = -2.542432795193405
$ws.Range("H11").Value = 12.16191994541689

# Row 12
$ws.Range("D12").Value = 113271.8445632598
$ws.Range("E12").Value = -0.09227009372118057
$ws.Range("F12").Value = 0.1970609302702817
$ws.Range("G12").Value = -1.104622018013555
$ws.Range("H12").Value = 8.165384038332004

# Row 15
$ws.Range("D15").Value = 113304.0261153597
$ws.Range("E15").Value = -0.07971892629489022
$ws.Range("F15").Value = 0.1711120050011246
$ws.Range("G15").Value = -0.9924053257321986
$ws.Range("H15").Value = 9.64387354538392
